$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values (shared strings): "batman" -> "username", "wsad" -> "adadsd"
$ws.Range("A1").Value = "username"
$ws.Range("B2").Value = "adadsd"

# Row 1 header (A1:B1) takes on the existing blue hyperlink-style font
$ws.Range("A1:B1").Font.Color = $ws.Range("A3").Font.Color

# Column A rows 3-5 (hyperlinked emails) take on a new bold font style,
# built from a default-styled cell so no stray color is inherited
$ws.Range("B3").Font.Bold = $true
$ws.Range("B3").Copy()
$ws.Range("A3:A5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B3").Font.Bold = $false
